# Apply the F-column ("想去人数") updates to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 310
$ws1.Range("F3").Value = 47
$ws1.Range("F5").Value = 4545
$ws1.Range("F9").Value = 704
$ws1.Range("F10").Value = 187

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 310
$ws4.Range("F3").Value = 47
$ws4.Range("F5").Value = 4545
$ws4.Range("F9").Value = 704
$ws4.Range("F11").Value = 187
